# Finalized Experiments with Participant Generation
# Rename task-order sheets and update their generated filename/timestamp cells.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16502912293889132"
$ws1.Range("B2").Value = "go_stims-16502912293489127.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912293719115.csv"
$ws1.Range("B4").Value = "go_stims-16502912293729305.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912293869445.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16502912331872728"
$ws2.Range("B2").Value = "OB-1650291231767992.csv"
$ws2.Range("B3").Value = "ZB-match_1-1650291229563951.csv"
$ws2.Range("B4").Value = "TB-16502912324832704.csv"
$ws2.Range("B5").Value = "ZB-match_0-16502912303896868.csv"
$ws2.Range("B6").Value = "OB-16502912304099886.csv"
$ws2.Range("B7").Value = "TB-16502912319582744.csv"
$ws2.Range("B8").Value = "OB-16502912309570258.csv"
$ws2.Range("B9").Value = "TB-16502912331702704.csv"
$ws2.Range("B10").Value = "ZB-match_8-16502912302496889.csv"

# --- Sheet 3: RS_TO (name only, content unchanged) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16502912331892745"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16502912332362697"
$ws4.Range("B2").Value = "MM_stims-16502912332032778.csv"
$ws4.Range("B3").Value = "ZM_stims-165029123319127.csv"
$ws4.Range("B4").Value = "MM_stims-1650291233219307.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912332042708.csv"
$ws4.Range("B6").Value = "MM_stims-16502912332362697.csv"
$ws4.Range("B7").Value = "ZM_stims-1650291233219307.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16502912332972677"
$ws5.Range("B2").Value = "SAT_stims-16502912332402678.csv"
$ws5.Range("B3").Value = "vSAT_stims-165029123326727.csv"
$ws5.Range("B4").Value = "SAT_stims-16502912332522688.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502912332823017.csv"
